$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Unite" column (column C, "EUR/MWh") entirely.
# Columns D (16-jun) and E (17-jun) shift left to become C and D.
$ws.Range("C:C").Delete()

# Remove the trailing summary rows (FR DAY PEAK, FR DAY BASE, JOURNEE ECO...)
$ws.Range("26:28").Delete()

Write-Output "done"
